# Add a new "keterangan" column to the product upload/download template.
#
# The template's data table occupies columns A:H (barcode, nama, satuan,
# stok, harga_beli, harga_ecer, harga_grosir, min_beli_grosir) in row 1,
# with the "catatan:" notes block living off to the right in column K.
# This change inserts a brand-new column I, headed "keterangan", right
# after the existing data columns. Inserting (rather than just writing
# into a blank column) correctly shifts the pre-existing notes block
# from column K to column L, which is the behavior captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (9); this pushes the old column I (which was
# empty) and everything to its right (including the notes in K) one
# column to the right, so the notes end up in column L.
$ws.Columns.Item(9).Insert()

# Give the new column its header text, using the same header styling
# as the rest of row 1 (bold + border), which the sheet already applies
# automatically to inserted cells that inherit the row's formatting.
$ws.Range("I1").Value2 = "keterangan"

# Match the width used for the other note/data columns (same visual
# sizing as column H).
$ws.Columns.Item(9).ColumnWidth = 13.83

# Reflect the new active selection noted in the saved workbook.
$ws.Range("L13").Select()
